$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 2, 1.02),
    @(2, 3, 1.020739023565577),
    @(2, 4, 1.023383788850451),
    @(2, 5, 1.021710124372327),
    @(2, 6, 1.032008871629843),
    @(2, 9, 1.030658878237871),
    @(2, 10, 1.025934434659937),
    @(2, 11, 1.026215386724597),
    @(2, 12, 1.024546651824302),
    @(2, 13, 1.034815340958775),
    @(2, 14, 1.012723153298275),
    @(3, 2, 1.02),
    @(3, 3, 1.021647206146306),
    @(3, 4, 1.024204234389573),
    @(3, 5, 1.022478920721165),
    @(3, 6, 1.033150557738757),
    @(3, 9, 1.030815290256242),
    @(3, 10, 1.026480149982627),
    @(3, 11, 1.026842662581355),
    @(3, 12, 1.025122062327767),
    @(3, 13, 1.035764818409987),
    @(3, 14, 1.012903942161363),
    @(4, 2, 1.02),
    @(4, 3, 1.022234911755387),
    @(4, 4, 1.024735450234839),
    @(4, 5, 1.022976834058881),
    @(4, 6, 1.033889239151937),
    @(4, 9, 1.030914482627918),
    @(4, 10, 1.026832713393403),
    @(4, 11, 1.027248239625685),
    @(4, 12, 1.025494185369543),
    @(4, 13, 1.036378553493329),
    @(4, 14, 1.01302071631921),
    @(5, 2, 1.02),
    @(5, 3, 1.022481994522112),
    @(5, 4, 1.02495885184943),
    @(5, 5, 1.023186263498707),
    @(5, 6, 1.034199764963797),
    @(5, 9, 1.030955699683256),
    @(5, 10, 1.026980798345791),
    @(5, 11, 1.027418668495763),
    @(5, 12, 1.025650575666081),
    @(5, 13, 1.036636414159825),
    @(5, 14, 1.013069757998456),
    @(6, 2, 1.02),
    @(6, 3, 1.022523481446248),
    @(6, 4, 1.024996366530801),
    @(6, 5, 1.023221433842818),
    @(6, 6, 1.034251902681533),
    @(6, 9, 1.03096259184208),
    @(6, 10, 1.027005654655028),
    @(6, 11, 1.027447279785455),
    @(6, 12, 1.025676831313048),
    @(6, 13, 1.036679701087449),
    @(6, 14, 1.013077989360957),
    @(7, 2, 1.02),
    @(7, 3, 1.022238213243929),
    @(7, 4, 1.024738435031294),
    @(7, 5, 1.022979632047703),
    @(7, 6, 1.033893388475463),
    @(7, 9, 1.030915035272525),
    @(7, 10, 1.026834692634267),
    @(7, 11, 1.027250517203299),
    @(7, 12, 1.025496275262858),
    @(7, 13, 1.03638199964393),
    @(7, 14, 1.013021371814073),
    @(8, 2, 1.02),
    @(8, 3, 1.02104593712012),
    @(8, 4, 1.023660992662177),
    @(8, 5, 1.021969848645305),
    @(8, 6, 1.032394723375693),
    @(8, 9, 1.0307121556019),
    @(8, 10, 1.026118975283159),
    @(8, 11, 1.026427441657381),
    @(8, 12, 1.024741156615038),
    @(8, 13, 1.0351363537312),
    @(8, 14, 1.012784294694592),
    @(9, 2, 1.02),
    @(9, 3, 1.018945418717983),
    @(9, 4, 1.021765001447214),
    @(9, 5, 1.020193986672613),
    @(9, 6, 1.029753382813528),
    @(9, 9, 1.030339237602121),
    @(9, 10, 1.024853606455448),
    @(9, 11, 1.024974725760659),
    @(9, 12, 1.023409002741984),
    @(9, 13, 1.032936473421835),
    @(9, 14, 1.012364952232903),
    @(10, 2, 1.02),
    @(10, 3, 1.017545404543967),
    @(10, 4, 1.020502823406391),
    @(10, 5, 1.019012506527999),
    @(10, 6, 1.027992152282702),
    @(10, 9, 1.030080294813624),
    @(10, 10, 1.024007261340852),
    @(10, 11, 1.024004720957709),
    @(10, 12, 1.022519917994995),
    @(10, 13, 1.031466608074518),
    @(10, 14, 1.012084344137782),
    @(11, 2, 1.02),
    @(11, 3, 1.016939270572932),
    @(11, 4, 1.019956731000265),
    @(11, 5, 1.018501503003402),
    @(11, 6, 1.027229438863479),
    @(11, 9, 1.029965726518161),
    @(11, 10, 1.023640137937234),
    @(11, 11, 1.023584346681489),
    @(11, 12, 1.022134712117656),
    @(11, 13, 1.03082936380731),
    @(11, 14, 1.011962593049338),
    @(12, 2, 1.02),
    @(12, 3, 1.016714138303317),
    @(12, 4, 1.019753954989282),
    @(12, 5, 1.018311782592889),
    @(12, 6, 1.026946119519924),
    @(12, 9, 1.02992280402919),
    @(12, 10, 1.02350367494797),
    @(12, 11, 1.023428148321015),
    @(12, 12, 1.021991596309441),
    @(12, 13, 1.030592545262173),
    @(12, 14, 1.011917332591101),
    @(13, 2, 1.02),
    @(13, 3, 1.01676242934642),
    @(13, 4, 1.019797448081974),
    @(13, 5, 1.018352474217068),
    @(13, 6, 1.027006893101298),
    @(13, 9, 1.029932027635565),
    @(13, 10, 1.023532951101144),
    @(13, 11, 1.023461655768449),
    @(13, 12, 1.022022296625322),
    @(13, 13, 1.030643348968558),
    @(13, 14, 1.011927042768777),
    @(14, 2, 1.02),
    @(14, 3, 1.016920660802921),
    @(14, 4, 1.019939968098539),
    @(14, 5, 1.018485818831932),
    @(14, 6, 1.027206019877477),
    @(14, 9, 1.029962186009914),
    @(14, 10, 1.023628859838927),
    @(14, 11, 1.023571436345362),
    @(14, 12, 1.022122882792805),
    @(14, 13, 1.030809790712753),
    @(14, 14, 1.011958852552536),
    @(15, 2, 1.02),
    @(15, 3, 1.017018154150922),
    @(15, 4, 1.020027788297539),
    @(15, 5, 1.018567988684032),
    @(15, 6, 1.0273287066439),
    @(15, 9, 1.029980719000647),
    @(15, 10, 1.023687939543664),
    @(15, 11, 1.023639068831427),
    @(15, 12, 1.022184852875865),
    @(15, 13, 1.030912325373624),
    @(15, 14, 1.011978446757214),
    @(16, 2, 1.02),
    @(16, 3, 1.017585633429628),
    @(16, 4, 1.020539075101589),
    @(16, 5, 1.019046432574977),
    @(16, 6, 1.02804276916437),
    @(16, 9, 1.030087846878285),
    @(16, 10, 1.024031612469254),
    @(16, 11, 1.024032612387638),
    @(16, 12, 1.022545478141263),
    @(16, 13, 1.0315088833809),
    @(16, 14, 1.012092419197751),
    @(17, 2, 1.02),
    @(17, 3, 1.017941620439646),
    @(17, 4, 1.020859910097474),
    @(17, 5, 1.019346705258806),
    @(17, 6, 1.0284906577207),
    @(17, 9, 1.03015439130645),
    @(17, 10, 1.024247015786843),
    @(17, 11, 1.024279377185448),
    @(17, 12, 1.022771628823412),
    @(17, 13, 1.031882878839648),
    @(17, 14, 1.012163845443768),
    @(18, 2, 1.02),
    @(18, 3, 1.018149269451126),
    @(18, 4, 1.021047090009902),
    @(18, 5, 1.01952190562021),
    @(18, 6, 1.028751894970177),
    @(18, 9, 1.030192969666182),
    @(18, 10, 1.0243725940983),
    @(18, 11, 1.024423276533532),
    @(18, 12, 1.022903516800243),
    @(18, 13, 1.032100948530009),
    @(18, 14, 1.012205483418128),
    @(19, 2, 1.02),
    @(19, 3, 1.018220073703163),
    @(19, 4, 1.021110920663007),
    @(19, 5, 1.019581653945122),
    @(19, 6, 1.028840968675613),
    @(19, 9, 1.030206083864917),
    @(19, 10, 1.024415402383971),
    @(19, 11, 1.024472336635534),
    @(19, 12, 1.022948483454499),
    @(19, 13, 1.03217529184817),
    @(19, 14, 1.012219676855437),
    @(20, 2, 1.02),
    @(20, 3, 1.017903425589785),
    @(20, 4, 1.020825483157012),
    @(20, 5, 1.019314483001151),
    @(20, 6, 1.028442604405217),
    @(20, 9, 1.030147276113421),
    @(20, 10, 1.024223911539922),
    @(20, 11, 1.024252905208487),
    @(20, 12, 1.022747367253776),
    @(20, 13, 1.031842760479214),
    @(20, 14, 1.012156184537818),
    @(21, 2, 1.02),
    @(21, 3, 1.016874065222209),
    @(21, 4, 1.019897997647172),
    @(21, 5, 1.018446549699056),
    @(21, 6, 1.027147382383021),
    @(21, 9, 1.029953315239591),
    @(21, 10, 1.023600619778481),
    @(21, 11, 1.023539110148179),
    @(21, 12, 1.022093263583389),
    @(21, 13, 1.030760781015579),
    @(21, 14, 1.011949486374566),
    @(22, 2, 1.02),
    @(22, 3, 1.016226940313787),
    @(22, 4, 1.019315238995363),
    @(22, 5, 1.017901361362332),
    @(22, 6, 1.026332945961588),
    @(22, 9, 1.029829242844769),
    @(22, 10, 1.023208170577465),
    @(22, 11, 1.023090014721437),
    @(22, 12, 1.021681810989399),
    @(22, 13, 1.030079816869979),
    @(22, 14, 1.011819314928803),
    @(23, 2, 1.02),
    @(23, 3, 1.016569986166092),
    @(23, 4, 1.019624133222139),
    @(23, 5, 1.018190326747835),
    @(23, 6, 1.026764701685495),
    @(23, 9, 1.029895216891838),
    @(23, 10, 1.023416268261326),
    @(23, 11, 1.023328117262179),
    @(23, 12, 1.021899947693054),
    @(23, 13, 1.030440873384379),
    @(23, 14, 1.0118883412923),
    @(24, 2, 1.02),
    @(24, 3, 1.017920684174205),
    @(24, 4, 1.020841039078922),
    @(24, 5, 1.019329042678262),
    @(24, 6, 1.028464317656287),
    @(24, 9, 1.030150491891721),
    @(24, 10, 1.024234351548656),
    @(24, 11, 1.024264866861678),
    @(24, 12, 1.022758330080447),
    @(24, 13, 1.031860888472616),
    @(24, 14, 1.012159646244842),
    @(25, 2, 1.02),
    @(25, 3, 1.019488398012738),
    @(25, 4, 1.022254845479894),
    @(25, 5, 1.020652666809058),
    @(25, 6, 1.030436292012391),
    @(25, 9, 1.030437469287571),
    @(25, 10, 1.025181225538731),
    @(25, 11, 1.025350560777216),
    @(25, 12, 1.023753573356644),
    @(25, 13, 1.033505773868054),
    @(25, 14, 1.012473548020191)
)

foreach ($row in $data) {
    $r = $row[0]
    $c = $row[1]
    $v = $row[2]
    $ws.Cells.Item($r, $c).Value = $v
}
